$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "TreatmentTab" query cell (B5) wraps its REPLACE() call in a redundant
# CONCAT(...) - remove the superfluous CONCAT wrapper while keeping the rest
# of the SQL query text identical.
$treatmentCell = $ws.Range("B5")
$oldQuery = $treatmentCell.Value()
$newQuery = $oldQuery.Replace(
    "CONCAT(REPLACE(trt.treatment_agent, ';', ', ')) AS ""Treatment Agent"",",
    "REPLACE(trt.treatment_agent, ';', ', ') AS ""Treatment Agent"","
)
$treatmentCell.Value = $newQuery
